$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value = "0.9999992176520535" -as [double]
$ws.Range("E3").Value = "0.9999992176520535" -as [double]

$ws.Range("D4").Value = "6.125234272457923E-10" -as [double]
$ws.Range("E4").Value = "6.125234272457923E-10" -as [double]

$ws.Range("D5").Value = "9.515590701970698E-08" -as [double]
$ws.Range("E5").Value = "9.515590701970698E-08" -as [double]

$ws.Range("D6").Value = "1.312496810941614E-58" -as [double]
$ws.Range("E6").Value = "1.312496810941614E-58" -as [double]

$ws.Range("D10").Value = "1.161779723084769E-15" -as [double]
$ws.Range("E10").Value = "0.9999999999999989" -as [double]

$ws.Range("D11").Value = "1.389074516416762E-05" -as [double]
$ws.Range("E11").Value = "0.9999861092548359" -as [double]
$ws.Range("F11").Value = "126.1665496826172" -as [double]
